$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 8 ("3ª Dose" / "3ª dose"),
# shifting rows 8-17 down to 9-18, then fill the new row 8 with
# ("Única", "Dose única").
$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = "Única"
$ws.Range("B8").Value = "Dose única"

# Match the page setup that appears in the saved file (paper size +
# portrait orientation).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
